$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCols = "A","B","C","D","E","F","G","K","L"
foreach ($col in $textCols) {
    foreach ($r in 2..5) {
        $ws.Range("$col$r").NumberFormat = "@"
    }
}

# Row 2
$ws.Range("A2").Value = "281474993058530-1743637850190"
$ws.Range("B2").Value = "Mobile Usage"
$ws.Range("C2").Value = "2025-04-02T17:50:50.190"
$ws.Range("D2").Value = "281474993058530"
$ws.Range("E2").Value = "154"
$ws.Range("F2").Value = "51834043"
$ws.Range("G2").Value = "MIGUEL ÁNGEL GUIZAR"
$ws.Range("K2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743637847690/PaZ2GKc6aN-camera-video-segment-driver-1743637850190.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T170139Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2001%3A01%3A39%20GMT&X-Amz-Signature=7c1e4dad916ec5a203878ce3677cd5d5bfff8b4b981865b9c8f90c490cc28b1a"
$ws.Range("L2").Value = "No video URL"
$ws.Range("H2").Value = 20.934953509
$ws.Range("I2").Value = -105.32674498
$ws.Range("J2").Value = 0

# Row 3
$ws.Range("A3").Value = "281474993058530-1743634724101"
$ws.Range("B3").Value = "Mobile Usage"
$ws.Range("C3").Value = "2025-04-02T16:58:44.101"
$ws.Range("D3").Value = "281474993058530"
$ws.Range("E3").Value = "154"
$ws.Range("F3").Value = "51834043"
$ws.Range("G3").Value = "MIGUEL ÁNGEL GUIZAR"
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474993058530/1743634721601/5dyCDJLTC3-camera-video-segment-driver-1743634724101.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T170139Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2001%3A01%3A39%20GMT&X-Amz-Signature=1e7d7d173fa71928a8ed1ec6ca3281e583479767f7cab72aa86e5eae7e6c0421"
$ws.Range("L3").Value = "No video URL"
$ws.Range("H3").Value = 20.707228419
$ws.Range("I3").Value = -105.27417089
$ws.Range("J3").Value = 0

# Row 4
$ws.Range("A4").Value = "281474991395097-1743627062321"
$ws.Range("B4").Value = "Harsh Brake"
$ws.Range("C4").Value = "2025-04-02T14:51:02.321"
$ws.Range("D4").Value = "281474991395097"
$ws.Range("E4").Value = "125"
$ws.Range("F4").Value = "51834055"
$ws.Range("G4").Value = "DAVID SERRANO"
$ws.Range("K4").Value = "No video URL"
$ws.Range("L4").Value = "No video URL"
$ws.Range("H4").Value = 20.65682646
$ws.Range("I4").Value = -103.37152569
$ws.Range("J4").Value = 0.7331250309944153

# Row 5
$ws.Range("A5").Value = "281474991109864-1743618609572"
$ws.Range("B5").Value = "Forward Collision Warning"
$ws.Range("C5").Value = "2025-04-02T12:30:09.572"
$ws.Range("D5").Value = "281474991109864"
$ws.Range("E5").Value = "138"
$ws.Range("F5").Value = "51833957"
$ws.Range("G5").Value = "FERNANDO ORNELAS"
$ws.Range("K5").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1743618603072/lsfY7Elc8s-camera-video-segment-driver-1743618608072.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T170139Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2001%3A01%3A39%20GMT&X-Amz-Signature=0cf6312809aa07fae51626075980a4434e80f5f6c963464f9f0422b638746be9"
$ws.Range("L5").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743618603072/krpw3PxSiC-camera-video-segment-1743618608072.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK65G53T%2F20250403%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250403T170139Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCICgh8znk17t3TA%2BoKV6Gs6dr67pjsww%2B74TUfK2AKgf7AiEAz%2Bsh73bZQ2b5CCsKg0Vjdzko8tsT%2FHUnEee9To43qz0q5gMI8f%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDCtu7H5t94rCmyRPWSq6A40NpyKaxkLrPb21f0WvTGT%2BIGiyn7Vbahoust2%2FGVAklzVYhIFiLefEp7yFI%2B7oJscwdjq7YKspHYMijhn7XRWKB4q464bVv7AXR5Uk5%2FmcOBPpkQAl0Be5gWWOYqROqpt%2FfDBW4shLaQ%2B2Djfz8obOjliQROGSRIdhfFY2BupuWf1VQ8prpoAvXUduCB3rxfGPNwH89vgp9nigybAXOLa25dLV1oq%2BzBxb39qnxbzdW8VQMOeQCwvTX%2Fq0vDxuDeN0iQREY4vwy%2BL13iBqN%2BkPYF9Jctlv43thbAG6o58%2Fh%2BOQQtWOGAykxPLN4qYJNqUBr%2BXgd1UBqMrKKgyIrmtMyVDf1zZYCtD7TuGL2AIszRTuVQY5tEIt2kxWCKFfKSymj8xk6kI%2ByUwaDrztK2xDDGJ%2Bqpe7Lcn5m6MkSrqab%2BvdDQEk6vYmJNdixVjZT6Q7SzfLXAkUsHZGE76ZlcPr8ko7UkStEK06fM%2BCql0r%2BD%2FXGcC0jvV5hXcN3QACFNyeMargdFsewTxlm%2FpWJgfI7X4u0LO8goE%2F67DaL8MqGpWCVkbvMPrJVZTdQIEwip1Af0893leelaswhti6vwY6pQHL%2FVpWyuVSfKKmSWKrk0yqWCTU%2FIZQrKI6CNfDOMGVTVeh%2BeQlMGeneOMBX83X3%2FtDD9GAWK2Grz19JgJk9GJufgKV75gfSafxWlhJ2gBKoWQZYRuoK%2Fo8ietECWT19yHU6Xs%2BbtZFlDaa346Gp%2BkCKvs3RDD%2FQ%2Bs%2BDQYMcujj%2B3UiGHSGthFgBpaAKoLr05Qg0%2BoOLYJ3JHH1JZtQYDvcgCTG6BY%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2004%20Apr%202025%2001%3A01%3A39%20GMT&X-Amz-Signature=16662a96b81474b060f93d28bf015cf7b27f8a63364c4eb8874c1231560a90a1"
$ws.Range("H5").Value = 20.70763009
$ws.Range("I5").Value = -103.45866694
$ws.Range("J5").Value = 0

# Reset style of text cells to Normal (default, no explicit style) while preserving text type
foreach ($col in $textCols) {
    foreach ($r in 2..5) {
        $ws.Range("$col$r").Style = "Normal"
    }
}
